# Append a new client row (row 7) to the "Klijenti" sheet, matching the
# existing columns: fullName, phone, city, address, email.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "Marko Petrović"
$ws.Range("B7").Value = "069/987-654"
$ws.Range("C7").Value = "Kotor"
$ws.Range("D7").Value = ""
$ws.Range("E7").Value = ""
